$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing row 2 (ECs -> Fgf9 -> Fgfr2 -> ECs) ---
$ws.Cells.Item(2, 7).Value = 3.710887666666667
$ws.Cells.Item(2, 8).Value = 11.132663
$ws.Cells.Item(2, 11).Value = 3
$ws.Cells.Item(2, 12).Value = 1
$ws.Cells.Item(2, 13).Value = 0.714474
$ws.Cells.Item(2, 14).Value = 2.143422
$ws.Cells.Item(2, 15).Value = 0.138796410342318
$ws.Cells.Item(2, 16).Value = 0.138796410342318
$ws.Cells.Item(2, 17).Value = 2.651332754754
$ws.Cells.Item(2, 18).Value = 23.861994792786
$ws.Cells.Item(2, 19).Value = 0.138796410342318
$ws.Cells.Item(2, 20).Value = 0.138796410342318

# --- Update existing row 3 (ECs -> Fgf9 -> Fgfr2 -> FAPs) ---
$ws.Cells.Item(3, 7).Value = 3.710887666666667
$ws.Cells.Item(3, 8).Value = 11.132663
$ws.Cells.Item(3, 15).Value = 0.8044215857867821
$ws.Cells.Item(3, 16).Value = 0.8044215857867821
$ws.Cells.Item(3, 17).Value = 15.366314544933
$ws.Cells.Item(3, 18).Value = 138.296830904397
$ws.Cells.Item(3, 19).Value = 0.8044215857867821
$ws.Cells.Item(3, 20).Value = 0.8044215857867821

# --- Update existing row 4 (ECs -> Fgf9 -> Fgfr2 -> MuSCs) ---
$ws.Cells.Item(4, 7).Value = 3.710887666666667
$ws.Cells.Item(4, 8).Value = 11.132663
$ws.Cells.Item(4, 13).Value = 0.2847646666666667
$ws.Cells.Item(4, 14).Value = 0.8542940000000001
$ws.Cells.Item(4, 15).Value = 0.05531945672713084
$ws.Cells.Item(4, 16).Value = 0.05531945672713083
$ws.Cells.Item(4, 17).Value = 1.056729689435778
$ws.Cells.Item(4, 18).Value = 9.510567204922001
$ws.Cells.Item(4, 19).Value = 0.05531945672713084
$ws.Cells.Item(4, 20).Value = 0.05531945672713083

# --- Add new row 5 (ECs -> Fgf9 -> Fgfr2 -> Resolving-Mac) ---
$ws.Cells.Item(5, 1).Value = "ECs"
$ws.Cells.Item(5, 2).Value = "Fgf9"
$ws.Cells.Item(5, 3).Value = "Fgfr2"
$ws.Cells.Item(5, 4).Value = "Resolving-Mac"
$ws.Cells.Item(5, 5).Value = 3
$ws.Cells.Item(5, 6).Value = 1
$ws.Cells.Item(5, 7).Value = 3.710887666666667
$ws.Cells.Item(5, 8).Value = 11.132663
$ws.Cells.Item(5, 9).Value = 1
$ws.Cells.Item(5, 10).Value = 1
$ws.Cells.Item(5, 11).Value = 2
$ws.Cells.Item(5, 12).Value = 0.6666666666666666
$ws.Cells.Item(5, 13).Value = 0.007528666666666667
$ws.Cells.Item(5, 14).Value = 0.022586
$ws.Cells.Item(5, 15).Value = 0.00146254714376898
$ws.Cells.Item(5, 16).Value = 0.00146254714376898
$ws.Cells.Item(5, 17).Value = 0.02793803627977778
$ws.Cells.Item(5, 18).Value = 0.251442326518
$ws.Cells.Item(5, 19).Value = 0.00146254714376898
$ws.Cells.Item(5, 20).Value = 0.00146254714376898
